$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2072368421052632
$ws.Range("C2").Value = 0.5131578947368421
$ws.Range("J2").Value = 0.01973684210526316
$ws.Range("P2").Value = 0.1644736842105263
$ws.Range("S2").Value = 0.09539473684210527
$ws.Range("B3").Value = 0.006329113924050633
$ws.Range("C3").Value = 0.0189873417721519
$ws.Range("J3").Value = 0.05063291139240506
$ws.Range("P3").Value = 0.7025316455696202
$ws.Range("S3").Value = 0.2215189873417721
$ws.Range("P4").Value = 0.7368421052631579
$ws.Range("S4").Value = 0.2631578947368421
$ws.Range("B6").Value = 0.07692307692307693
$ws.Range("D6").Value = 0.008547008547008548
$ws.Range("F6").Value = 0.05982905982905983
$ws.Range("J6").Value = 0.1837606837606838
$ws.Range("O6").Value = 0.01282051282051282
$ws.Range("Q6").Value = 0.1666666666666667
$ws.Range("R6").Value = 0.07692307692307693
$ws.Range("S6").Value = 0.4145299145299146
$ws.Range("B7").Value = 0.141025641025641
$ws.Range("D7").Value = 0.02564102564102564
$ws.Range("F7").Value = 0.05128205128205128
$ws.Range("J7").Value = 0.1282051282051282
$ws.Range("O7").Value = 0.02564102564102564
$ws.Range("Q7").Value = 0.1858974358974359
$ws.Range("R7").Value = 0.0641025641025641
$ws.Range("S7").Value = 0.3782051282051282
$ws.Range("B8").Value = 0.0825
$ws.Range("D8").Value = 0.025
$ws.Range("E8").Value = 0.0025
$ws.Range("F8").Value = 0.055
$ws.Range("J8").Value = 0.0925
$ws.Range("O8").Value = 0.015
$ws.Range("Q8").Value = 0.205
$ws.Range("R8").Value = 0.1075
$ws.Range("S8").Value = 0.415
$ws.Range("B9").Value = 0.1204188481675393
$ws.Range("D9").Value = 0.03141361256544502
$ws.Range("E9").Value = 0.005235602094240838
$ws.Range("F9").Value = 0.07329842931937172
$ws.Range("J9").Value = 0.03664921465968586
$ws.Range("O9").Value = 0.01570680628272251
$ws.Range("Q9").Value = 0.2303664921465969
$ws.Range("R9").Value = 0.1099476439790576
$ws.Range("S9").Value = 0.3769633507853403
$ws.Range("B10").Value = 0.1118991331757289
$ws.Range("D10").Value = 0.01339637509850276
$ws.Range("F10").Value = 0.07880220646178093
$ws.Range("J10").Value = 0.1142631993695823
$ws.Range("O10").Value = 0.01654846335697399
$ws.Range("Q10").Value = 0.2214342001576044
$ws.Range("R10").Value = 0.08825847123719464
$ws.Range("S10").Value = 0.355397951142632
$ws.Range("G11").Value = 0.1403508771929824
$ws.Range("J11").Value = 0.143859649122807
$ws.Range("K11").Value = 0.2315789473684211
$ws.Range("L11").Value = 0.4666666666666667
$ws.Range("S11").Value = 0.01754385964912281
$ws.Range("G12").Value = 0.697841726618705
$ws.Range("J12").Value = 0.1726618705035971
$ws.Range("K12").Value = 0.007194244604316547
$ws.Range("L12").Value = 0.04316546762589928
$ws.Range("S12").Value = 0.07913669064748201
$ws.Range("G13").Value = 0.6410256410256411
$ws.Range("J13").Value = 0.2307692307692308
$ws.Range("S13").Value = 0.1282051282051282
$ws.Range("F15").Value = 0.03017241379310345
$ws.Range("H15").Value = 0.125
$ws.Range("I15").Value = 0.04310344827586207
$ws.Range("J15").Value = 0.4439655172413793
$ws.Range("K15").Value = 0.06896551724137931
$ws.Range("M15").Value = 0.008620689655172414
$ws.Range("O15").Value = 0.06465517241379311
$ws.Range("S15").Value = 0.2155172413793103
$ws.Range("F16").Value = 0.03804347826086957
$ws.Range("H16").Value = 0.125
$ws.Range("I16").Value = 0.108695652173913
$ws.Range("J16").Value = 0.391304347826087
$ws.Range("K16").Value = 0.1467391304347826
$ws.Range("M16").Value = 0.0108695652173913
$ws.Range("N16").Value = 0.005434782608695652
$ws.Range("O16").Value = 0.08152173913043478
$ws.Range("S16").Value = 0.09239130434782608
$ws.Range("F17").Value = 0.0105708245243129
$ws.Range("H17").Value = 0.1374207188160677
$ws.Range("I17").Value = 0.09090909090909091
$ws.Range("J17").Value = 0.427061310782241
$ws.Range("K17").Value = 0.09725158562367865
$ws.Range("M17").Value = 0.01691331923890063
$ws.Range("N17").Value = 0.002114164904862579
$ws.Range("O17").Value = 0.08668076109936575
$ws.Range("S17").Value = 0.1310782241014799
$ws.Range("F18").Value = 0.004926108374384237
$ws.Range("H18").Value = 0.1773399014778325
$ws.Range("I18").Value = 0.06896551724137931
$ws.Range("J18").Value = 0.4088669950738916
$ws.Range("K18").Value = 0.1280788177339902
$ws.Range("M18").Value = 0.01477832512315271
$ws.Range("O18").Value = 0.07881773399014778
$ws.Range("S18").Value = 0.1182266009852217
$ws.Range("F19").Value = 0.01877551020408163
$ws.Range("H19").Value = 0.2048979591836735
$ws.Range("I19").Value = 0.08489795918367347
$ws.Range("J19").Value = 0.3893877551020408
$ws.Range("K19").Value = 0.08653061224489796
$ws.Range("M19").Value = 0.02040816326530612
$ws.Range("N19").Value = 0.00163265306122449
$ws.Range("O19").Value = 0.07428571428571429
$ws.Range("S19").Value = 0.1191836734693878
